$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.114.11'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = '3.313.80'
$ws.Range('E3').Value = '  -0.29%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '188.08'
$ws.Range('E5').Value = '  +4.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '558.89'
$ws.Range('E6').Value = '  +0.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.588'
$ws.Range('E8').Value = '  +0.27%  '
$ws.Range('D9').Value = '3.306.88'
$ws.Range('E9').Value = '  -0.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.184'
$ws.Range('E10').Value = '  +0.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.588'
$ws.Range('E11').Value = '  +1.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.63'
$ws.Range('E12').Value = '  +1.33%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000271'
$ws.Range('E13').Value = '  +3.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.67'
$ws.Range('E14').Value = '  +1.82%  '
$ws.Range('D15').Value = '3.858.51'
$ws.Range('E15').Value = '  -0.12%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '604.70'
$ws.Range('E16').Value = '  +1.07%  '
$ws.Range('D17').Value = '66.217.02'
$ws.Range('E17').Value = '  +0.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.02'
$ws.Range('E18').Value = '  +0.97%  '
$ws.Range('E19').Value = '  +1.31%  '
$ws.Range('D20').Value = '3.327.20'
$ws.Range('E20').Value = '  -0.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.14'
$ws.Range('E21').Value = '  -1.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.913'
$ws.Range('E22').Value = '  +1.73%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.50'
$ws.Range('E23').Value = '  +10.55%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.14'
$ws.Range('E24').Value = '  +1.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '100.27'
$ws.Range('E25').Value = '  +0.72%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.99'
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.06'
$ws.Range('E27').Value = '  +1.11%  '
$ws.Range('E28').Value = '  +5.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.67'
$ws.Range('E29').Value = '  +5.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.66'
$ws.Range('E30').Value = '  +0.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.44'
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.77'
$ws.Range('E32').Value = '  +9.18%  '
$ws.Range('E33').Value = '  +4.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '584.11'
$ws.Range('E34').Value = '  +10.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '11.15'
$ws.Range('E35').Value = '  +1.85%  '
$ws.Range('E36').Value = '  +2.13%  '
$ws.Range('D37').Value = '3.722.39'
$ws.Range('E37').Value = '  -1.34%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '57.26'
$ws.Range('E38').Value = '  -0.82%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  +0.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.53'
$ws.Range('E40').Value = '  +21.69%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.132'
$ws.Range('E41').Value = '  +6.31%  '
$ws.Range('B42').Value = 'PEPE'
$ws.Range('C42').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D42').Value = '0.0₃0726'
$ws.Range('E42').Value = '  +2.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '34.04'
$ws.Range('E43').Value = '  +7.91%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.29'
$ws.Range('E44').Value = '  -3.84%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.69'
$ws.Range('E45').Value = '  +1.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.344'
$ws.Range('E46').Value = '  +2.10%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0425'
$ws.Range('E47').Value = '  +3.57%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.33'
$ws.Range('E48').Value = '  +2.53%  '
$ws.Range('E49').Value = '  +1.19%  '
$ws.Range('E50').Value = '  +1.05%  '
$ws.Range('E51').Value = '  +0.09%  '
